$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.563.86"
$ws.Range("E2").Value = '  -2.35%  '
$ws.Range("D3").Value = "'2.002.41"
$ws.Range("E3").Value = '  -4.01%  '
$ws.Range("E4").Value = '  +0.81%  '
$ws.Range("D5").Value = "'329.57"
$ws.Range("E5").Value = '  -3.91%  '
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("D7").Value = "'0.5003"
$ws.Range("E7").Value = '  -4.34%  '
$ws.Range("D8").Value = "'0.4211"
$ws.Range("E8").Value = '  -4.42%  '
$ws.Range("D9").Value = "'54.09"
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").Value = "'0.09004"
$ws.Range("E10").Value = '  -3.66%  '
$ws.Range("D11").Value = "'1.118"
$ws.Range("E11").Value = '  -4.15%  '
$ws.Range("D12").Value = "'23.25"
$ws.Range("E12").Value = '  -6.02%  '
$ws.Range("D13").Value = "'2.030.05"
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").Value = "'8.039"
$ws.Range("E14").Value = '  -6.06%  '
$ws.Range("D15").Value = "'6.475"
$ws.Range("E15").Value = '  -5.96%  '
$ws.Range("D16").Value = "'1.013"
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").Value = "'94.35"
$ws.Range("D18").Value = "'0.00001112"
$ws.Range("E18").Value = '  -3.77%  '
$ws.Range("D19").Value = "'0.06681"
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").Value = "'19.66"
$ws.Range("E20").Value = '  -6.53%  '
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").Value = "'5.961"
$ws.Range("E22").Value = '  -5.58%  '
$ws.Range("D23").Value = "'29.603.86"
$ws.Range("E23").Value = '  -2.33%  '
$ws.Range("D24").Value = "'11.99"
$ws.Range("E24").Value = '  -4.00%  '
$ws.Range("D25").Value = "'2.305"
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = "'159.15"
$ws.Range("E26").Value = '  -1.94%  '
$ws.Range("D27").Value = "'20.71"
$ws.Range("E27").Value = '  -4.75%  '
$ws.Range("D28").Value = "'6.398"
$ws.Range("E28").Value = '  -4.15%  '
$ws.Range("D29").Value = "'2.297"
$ws.Range("E29").Value = '  -8.14%  '
$ws.Range("D30").Value = "'128.00"
$ws.Range("E30").Value = '  -3.65%  '
$ws.Range("D31").Value = "'1.055"
$ws.Range("E31").Value = '  -6.43%  '
$ws.Range("D32").Value = "'0.09956"
$ws.Range("E32").Value = '  -4.55%  '
$ws.Range("D33").Value = "'1.567"
$ws.Range("E33").Value = '  -5.70%  '
$ws.Range("D34").Value = "'5.834"
$ws.Range("E34").Value = '  -6.08%  '
$ws.Range("D35").Value = "'3.800"
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("D36").Value = "'0.02468"
$ws.Range("E36").Value = '  -5.73%  '
$ws.Range("D37").Value = "'9.266"
$ws.Range("E37").Value = '  -9.30%  '
$ws.Range("D38").Value = "'1.306"
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("D39").Value = "'0.06393"
$ws.Range("E39").Value = '  -6.16%  '
$ws.Range("D40").Value = "'0.6556"
$ws.Range("D41").Value = "'11.68"
$ws.Range("E41").Value = '  -6.26%  '
$ws.Range("D42").Value = "'0.2051"
$ws.Range("E42").Value = '  -6.87%  '
$ws.Range("E43").Value = '  +0.81%  '
$ws.Range("D44").Value = "'0.6365"
$ws.Range("E44").Value = '  -6.25%  '
$ws.Range("D45").Value = "'13.43"
$ws.Range("E45").Value = '  -6.00%  '
$ws.Range("D46").Value = "'2.189"
$ws.Range("E46").Value = '  -5.53%  '
$ws.Range("D47").Value = "'1.312"
$ws.Range("E47").Value = '  -4.05%  '
$ws.Range("D48").Value = "'3.510"
$ws.Range("E48").Value = '  -3.57%  '
$ws.Range("D49").Value = "'0.00000000340"
$ws.Range("E49").Value = '  -2.17%  '
$ws.Range("D50").Value = "'0.06986"
$ws.Range("E50").Value = '  -3.25%  '
$ws.Range("D51").Value = "'1.127"
$ws.Range("E51").Value = '  -6.31%  '
